$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header columns (row 1) to English snake_case field names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Normalize Spanish particles (de/del/la/las/los/el/y) to title case
#    in municipality/state names, plus a couple of one-off casing fixes
#    (GUANAJUATO -> Guanajuato, MonteMorelos -> Montemorelos).
$ws.Range("B6").Value = 'Pabellón De Arteaga'
$ws.Range("B7").Value = 'Rincón De Romos'
$ws.Range("B8").Value = 'San Francisco De Los Romo'
$ws.Range("B19").Value = 'Amatenango De La Frontera'
$ws.Range("B24").Value = 'Chiapa De Corzo'
$ws.Range("B31").Value = 'Mazapa De Madero'
$ws.Range("B32").Value = 'Montecristo De Guerero'
$ws.Range("B34").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B36").Value = 'San Cristóbal De Las Casas'
$ws.Range("B57").Value = 'Guadalupe Y Calvo'
$ws.Range("B59").Value = 'Hidalgo Del Parral'
$ws.Range("B73").Value = 'San Juan De Sabinas'
$ws.Range("A82").Value = 'Ciudad De México'
$ws.Range("B86").Value = 'Cuajimalpa De Morelos'
$ws.Range("B108").Value = 'Nombre De Dios'
$ws.Range("B110").Value = 'Pánuco De Coronado'
$ws.Range("B114").Value = 'San Juan De Guadalupe'
$ws.Range("B115").Value = 'San Juan Del Río'
$ws.Range("A122").Value = 'Estado De México'
$ws.Range("B122").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B124").Value = 'Almoloya De Alquisiras'
$ws.Range("B129").Value = 'Atizapán De Zaragoza'
$ws.Range("B131").Value = 'Chapa De Mota'
$ws.Range("B135").Value = 'Ecatepec De Morelos'
$ws.Range("B136").Value = 'Ixtapan De La Sal'
$ws.Range("B141").Value = 'Naucalpan De Juárez'
$ws.Range("B148").Value = 'San Felipe Del Progreso'
$ws.Range("B153").Value = 'Tenango Del Valle'
$ws.Range("B159").Value = 'Tlalnepantla De Baz'
$ws.Range("B163").Value = 'Valle De Bravo'
$ws.Range("B164").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B165").Value = 'Villa De Allende'
$ws.Range("A170").Value = 'Guanajuato'
$ws.Range("B173").Value = 'Apaseo El Alto'
$ws.Range("B174").Value = 'Apaseo El Grande'
$ws.Range("B181").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B184").Value = 'Jaral Del Progreso'
$ws.Range("B193").Value = 'San Diego De La Unión'
$ws.Range("B195").Value = 'San Francisco Del Rincón'
$ws.Range("B197").Value = 'San Luis De La Paz'
$ws.Range("B199").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B200").Value = 'Silao De La Victoria'
$ws.Range("B205").Value = 'Valle De Santiago'
$ws.Range("B210").Value = 'Acapulco De Juárez'
$ws.Range("B212").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B213").Value = 'Alcozauca De Guerero'
$ws.Range("B217").Value = 'Atoyac De Álvarez'
$ws.Range("B218").Value = 'Ayutla De Los Libres'
$ws.Range("B220").Value = 'Buenavista De Cuéllar'
$ws.Range("B221").Value = 'Chilapa De Álvarez'
$ws.Range("B222").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B225").Value = 'Coyuca De Catalán'
$ws.Range("B227").Value = 'Cutzamala De Pinzón'
$ws.Range("B231").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B232").Value = 'Iguala De La Independencia'
$ws.Range("B233").Value = 'Zihuatanejo De Azueta'
$ws.Range("B235").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B242").Value = 'Taxco De Alarcón'
$ws.Range("B244").Value = 'Técpan De Galeana'
$ws.Range("B246").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B254").Value = 'Agua Blanca De Iturbide'
$ws.Range("B258").Value = 'Atotonilco El Grande'
$ws.Range("B262").Value = 'Cuautepec De Hinojosa'
$ws.Range("B266").Value = 'Jacala De Ledezma'
$ws.Range("B270").Value = 'Mixquiahuala De Juárez'
$ws.Range("B272").Value = 'Nopala De Villagrán'
$ws.Range("B273").Value = 'Omitlán De Juárez'
$ws.Range("B274").Value = 'Pachuca De Soto'
$ws.Range("B276").Value = 'Progreso De Obregón'
$ws.Range("B277").Value = 'Santiago De Anaya'
$ws.Range("B281").Value = 'Tenango De Doria'
$ws.Range("B283").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B286").Value = 'Tula De Allende'
$ws.Range("B287").Value = 'Tulancingo De Bravo'
$ws.Range("B288").Value = 'Zacualtipán De Ángeles'
$ws.Range("B292").Value = 'Ahualulco De Mercado'
$ws.Range("B297").Value = 'Atemajac De Brizuela'
$ws.Range("B299").Value = 'Atotonilco El Alto'
$ws.Range("B300").Value = 'Autlán De Navarro'
$ws.Range("B307").Value = 'Huejuquilla El Alto'
$ws.Range("B308").Value = 'Ixtlahuacán De Los Membrillos'
$ws.Range("B310").Value = 'Lagos De Moreno'
$ws.Range("B314").Value = 'Ojuelos De Jalisco'
$ws.Range("B316").Value = 'San Diego De Alejandría'
$ws.Range("B319").Value = 'San Martín De Bolaños'
$ws.Range("B321").Value = 'San Miguel El Alto'
$ws.Range("B322").Value = 'Tamazula De Gordiano'
$ws.Range("B323").Value = 'Tepatitlán De Morelos'
$ws.Range("B328").Value = 'Unión De San Antonio'
$ws.Range("B329").Value = 'Valle De Guadalupe'
$ws.Range("B332").Value = 'Zapotlán El Grande'
$ws.Range("B376").Value = 'Tiquicheo De Nicolás Romero'
$ws.Range("B391").Value = 'Coatlán Del Río'
$ws.Range("B396").Value = 'Jonacatepec De Leandro Valle'
$ws.Range("B399").Value = 'Puente De Ixtla'
$ws.Range("B418").Value = 'Montemorelos'
$ws.Range("B420").Value = 'San Nicolás De Los Garza'
$ws.Range("B422").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B426").Value = 'Ayoquezco De Aldama'
$ws.Range("B431").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B432").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B434").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B437").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B438").Value = 'Oaxaca De Juárez'
$ws.Range("B439").Value = 'Ocotlán De Morelos'
$ws.Range("B441").Value = 'Putla Villa De Guerero'
$ws.Range("B443").Value = 'San Antonio De La Cal'
$ws.Range("B448").Value = 'San José Del Progreso'
$ws.Range("B480").Value = 'Santo Domingo De Morelos'
$ws.Range("B486").Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range("B487").Value = 'Tataltepec De Valdés'
$ws.Range("B488").Value = 'Tlacolula De Matamoros'
$ws.Range("B489").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B497").Value = 'Ayotoxco De Guerero'
$ws.Range("B505").Value = 'Cuetzalan Del Progreso'
$ws.Range("B510").Value = 'Huehuetlán El Chico'
$ws.Range("B512").Value = 'Ixcamilpa De Guerero'
$ws.Range("B514").Value = 'Izúcar De Matamoros'
$ws.Range("B530").Value = 'Tepexi De Rodríguez'
$ws.Range("B537").Value = 'Xochitlán De Vicente Suárez'
$ws.Range("B545").Value = 'Amealco De Bonfil'
$ws.Range("B550").Value = 'Jalpan De Serra'
$ws.Range("B551").Value = 'Landa De Matamoros'
$ws.Range("B553").Value = 'Pinal De Amoles'
$ws.Range("B555").Value = 'San Juan Del Río'
$ws.Range("B563").Value = 'Axtla De Terrazas'
$ws.Range("B572").Value = 'Mexquitic De Carmona'
$ws.Range("B577").Value = 'San Ciro De Acosta'
$ws.Range("B582").Value = 'Santa María Del Río'
$ws.Range("B583").Value = 'Soledad De Graciano Sánchez'
$ws.Range("B587").Value = 'Tanquián De Escobedo'
$ws.Range("B590").Value = 'Villa De Arista'
$ws.Range("B591").Value = 'Villa De Arriaga'
$ws.Range("B592").Value = 'Villa De Guadalupe'
$ws.Range("B593").Value = 'Villa De La Paz'
$ws.Range("B594").Value = 'Villa De Ramos'
$ws.Range("B595").Value = 'Villa De Reyes'
$ws.Range("B611").Value = 'Nacozari De García'
$ws.Range("B642").Value = 'Soto La Marina'
$ws.Range("B650").Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("B668").Value = 'Castillo De Teayo'
$ws.Range("B670").Value = 'Cazones De Herrera'
$ws.Range("B679").Value = 'Cosamaloapan De Carpio'
$ws.Range("B689").Value = 'Juchique De Ferrer'
$ws.Range("B691").Value = 'Lerdo De Tejada'
$ws.Range("B692").Value = 'Martínez De La Torre'
$ws.Range("B703").Value = 'Poza Rica De Hidalgo'
$ws.Range("B717").Value = 'Vega De Alatorre'
$ws.Range("B730").Value = 'Cañitas De Felipe Pescador'
$ws.Range("B731").Value = 'El Plateado De Joaquín Amaro'
$ws.Range("B751").Value = 'Tlaltenango De Sánchez Román'
$ws.Range("B753").Value = 'Villa De Cos'

# 3. Remove the trailing footnote/metadata rows (760-764) that are
#    no longer part of the clean dataset, shrinking the sheet to D758.
$ws.Range("A760:D764").EntireRow.Delete()

